$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 692.7273
$ws.Range("I6").Value = 124.44444
$ws.Range("J6").Value = 3250
$ws.Range("K6").Value = 373.33332
$ws.Range("L6").Value = 9750
$ws.Range("M6").Value = -261.33332
$ws.Range("N6").Value = -9974
$ws.Range("H44").Value = 40012.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 40012.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 40012.5
$ws.Range("N44").Value = -40936.5
$ws.Range("H96").Value = 1154.5454
$ws.Range("I96").Value = 300
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -17746
$ws.Range("H129").Value = 825.1724
$ws.Range("I129").Value = 356.57144
$ws.Range("J129").Value = 974.2727
$ws.Range("K129").Value = 1069.71432
$ws.Range("L129").Value = 2922.8181
$ws.Range("M129").Value = 3930.28568
$ws.Range("N129").Value = -12922.8181
$ws.Range("H132").Value = 58505.13
$ws.Range("I132").Value = 65494.824
$ws.Range("J132").Value = 4335
$ws.Range("K132").Value = 196484.472
$ws.Range("L132").Value = 13005
$ws.Range("M132").Value = -193954.472
$ws.Range("N132").Value = -18065
$ws.Range("H138").Value = 2117.04
$ws.Range("I138").Value = 991.36365
$ws.Range("J138").Value = 2671.4775
$ws.Range("K138").Value = 2974.09095
$ws.Range("L138").Value = 8014.4325
$ws.Range("M138").Value = 2165.90905
$ws.Range("N138").Value = -18294.4325

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 13854391
$ws.Range("I63").Value = 19789758
$ws.Range("J63").Value = 5200
$ws.Range("K63").Value = 19789758
$ws.Range("L63").Value = 5200
$ws.Range("M63").Value = -19789072
$ws.Range("N63").Value = -6572
$ws.Range("H66").Value = 13854391
$ws.Range("I66").Value = 19789758
$ws.Range("J66").Value = 5200
$ws.Range("K66").Value = 98948790
$ws.Range("L66").Value = 26000
$ws.Range("M66").Value = -98945358
$ws.Range("N66").Value = -32864
$ws.Range("H110").Value = 7170.3335
$ws.Range("I110").Value = 10005.5
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 10005.5
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = -7960.5
$ws.Range("N110").Value = -5590
$ws.Range("H122").Value = 2297.7778
$ws.Range("I122").Value = 1477.6842
$ws.Range("J122").Value = 4245.5
$ws.Range("K122").Value = 4433.0526
$ws.Range("L122").Value = 12736.5
$ws.Range("M122").Value = -1983.0526
$ws.Range("N122").Value = -17636.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2886.389
$ws.Range("I134").Value = 1744
$ws.Range("J134").Value = 5599.5625
$ws.Range("K134").Value = 5232
$ws.Range("L134").Value = 16798.6875
$ws.Range("M134").Value = -2697
$ws.Range("N134").Value = -21868.6875

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2690.6667
$ws.Range("I31").Value = 1230.2941
$ws.Range("J31").Value = 5611.4116
$ws.Range("K31").Value = 1230.2941
$ws.Range("L31").Value = 5611.4116
$ws.Range("M31").Value = -935.2941000000001
$ws.Range("N31").Value = -6201.4116
$ws.Range("H34").Value = 2690.6667
$ws.Range("I34").Value = 1230.2941
$ws.Range("J34").Value = 5611.4116
$ws.Range("K34").Value = 1230.2941
$ws.Range("L34").Value = 5611.4116
$ws.Range("M34").Value = -1028.2941
$ws.Range("N34").Value = -6015.4116
$ws.Range("H122").Value = 1970.4642
$ws.Range("I122").Value = 1411.32
$ws.Range("J122").Value = 6630
$ws.Range("K122").Value = 4233.96
$ws.Range("L122").Value = 19890
$ws.Range("M122").Value = -1783.96
$ws.Range("N122").Value = -24790
$ws.Range("H132").Value = 3292.3635
$ws.Range("I132").Value = 3247
$ws.Range("J132").Value = 3346.8
$ws.Range("K132").Value = 9741
$ws.Range("L132").Value = 10040.4
$ws.Range("M132").Value = -7211
$ws.Range("N132").Value = -15100.4
$ws.Range("H141").Value = 15710.143
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 15710.143
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 15710.143
$ws.Range("N141").Value = -26070.143

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5285.9
$ws.Range("I56").Value = 5285.9
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 5285.9
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -4755.9
$ws.Range("H113").Value = 521.475
$ws.Range("I113").Value = 530.6667
$ws.Range("J113").Value = 507.6875
$ws.Range("K113").Value = 1592.0001
$ws.Range("L113").Value = 1523.0625
$ws.Range("M113").Value = 577.9999
$ws.Range("N113").Value = -5863.0625
$ws.Range("H132").Value = 3406.2778
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 3950.9285
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 35558.3565
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -40618.3565

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1828.6666
$ws.Range("I13").Value = 1390
$ws.Range("J13").Value = 2048
$ws.Range("K13").Value = 1390
$ws.Range("L13").Value = 2048
$ws.Range("M13").Value = -1251
$ws.Range("N13").Value = -2326
$ws.Range("H122").Value = 5487.5
$ws.Range("I122").Value = 3983.3333
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 11949.9999
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -9499.999899999999
$ws.Range("N122").Value = -34900

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 28111.705
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 28111.705
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 28111.705
$ws.Range("N15").Value = -28687.705
$ws.Range("H113").Value = 8757.916999999999
$ws.Range("I113").Value = 11522.777
$ws.Range("J113").Value = 463.33334
$ws.Range("K113").Value = 34568.331
$ws.Range("L113").Value = 1390.00002
$ws.Range("M113").Value = -32398.331
$ws.Range("N113").Value = -5730.000019999999
$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 45000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -49900
$ws.Range("H126").Value = 675544.5600000001
$ws.Range("I126").Value = 3846.3333
$ws.Range("J126").Value = 1078563.5
$ws.Range("K126").Value = 11538.9999
$ws.Range("L126").Value = 3235690.5
$ws.Range("M126").Value = -9068.999899999999
$ws.Range("N126").Value = -3240630.5
$ws.Range("H132").Value = 5652985
$ws.Range("I132").Value = 4277.7856
$ws.Range("J132").Value = 10755043
$ws.Range("K132").Value = 12833.3568
$ws.Range("L132").Value = 32265129
$ws.Range("M132").Value = -10303.3568
$ws.Range("N132").Value = -32270189
